# Fix session persistence bug and update Finance Required bucket mappings
#
# 1) "buckets" sheet: rename the FIN_MAJOR/CORE bucket label from
#    "Core Required" to "Finance Required" and bump its needed_count
#    from 5 to 6 (the bucket now also requires every upper-division
#    FINA elective course, not just the 3 core FINA/ECON courses).
#
# 2) "bucket_course_map" sheet: the two ECON core rows (ECON 1103 /
#    ECON 1104) are replaced by the full list of upper-division FINA
#    elective course codes as CORE/required rows, while the existing
#    BUS_ELEC_4 elective rows for those same FINA courses are kept
#    (still optional / not required). This doubles-up the FINA course
#    list: once as a non-required business elective (BUS_ELEC_4) and
#    once as a required CORE course.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) buckets sheet
# ---------------------------------------------------------------
$buckets = $wb.Worksheets.Item("buckets")
$buckets.Cells.Item(2, 3).Value = "Finance Required"
$buckets.Cells.Item(2, 5).Value = 6

# ---------------------------------------------------------------
# 2) bucket_course_map sheet
# ---------------------------------------------------------------
$map = $wb.Worksheets.Item("bucket_course_map")

# The ordered list of upper-division FINA elective course codes that
# appear both as BUS_ELEC_4 electives (rows 50-70) and, after this
# change, as required CORE rows (rows 71-91).
$finaCodes = @(
    "FINA 3002",
    "FINA 4002",
    "FINA 4020",
    "FINA 4023",
    "FINA 4040",
    "FINA 4050",
    "FINA 4060",
    "FINA 4065",
    "FINA 4075",
    "FINA 4081",
    "FINA 4082",
    "FINA 4084",
    "FINA 4085",
    "FINA 4191",
    "FINA 4210",
    "FINA 4211",
    "FINA 4212",
    "FINA 4931",
    "FINA 4953",
    "FINA 4986",
    "FINA 4989"
)

# Row 49 (FIN_MAJOR / BUS_ELEC_4 / ECON 1104) is untouched.

# Rows 50-70: rewrite the former CORE/ECON + BUS_ELEC_4/FINA rows into a
# single contiguous BUS_ELEC_4 elective block covering every FINA code
# (not required, can double count).
$row = 50
for ($i = 0; $i -lt $finaCodes.Length; $i++) {
    $map.Cells.Item($row, 1).Value = "FIN_MAJOR"
    $map.Cells.Item($row, 2).Value = "BUS_ELEC_4"
    $map.Cells.Item($row, 3).Value = $finaCodes[$i]
    $map.Cells.Item($row, 4).Value = $false
    $map.Cells.Item($row, 5).Value = $true
    $row = $row + 1
}

# Rows 71-91: append the same FINA codes again, this time as required
# CORE rows.
for ($i = 0; $i -lt $finaCodes.Length; $i++) {
    $map.Cells.Item($row, 1).Value = "FIN_MAJOR"
    $map.Cells.Item($row, 2).Value = "CORE"
    $map.Cells.Item($row, 3).Value = $finaCodes[$i]
    $map.Cells.Item($row, 4).Value = $true
    $map.Cells.Item($row, 5).Value = $true
    $row = $row + 1
}
